# Updating plots and code
# - Remove Sheet2 (its unique data rows are no longer needed; deleting it also
#   prunes the now-unused shared strings that only Sheet2 referenced).
# - On Sheet1, reorder/relabel the summary rows and swap the pDNA row's
#   formulas from external-workbook references to literal hard-coded values.
# - Update the active selection on Sheet1.

$wb = $excel.ActiveWorkbook

# --- Remove Sheet2 ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: pDNA, now computed from literal values instead of the external
#     workbook reference ----------------------------------------------------
$ws.Range("A2").Value = "pDNA"
$ws.Range("B2").Formula = "=AVERAGE(-46.85,-52.23,-40.11)"
$ws.Range("C2").Formula = "=STDEV(-46.86,-52.24,-40.11)"

# --- Row 3: DIP S1 (G20:G22) ------------------------------------------------
$ws.Range("A3").Value = "DIP S1"
$ws.Range("B3").Formula = "=AVERAGE('[1]20231016_Zeta-PotentialData'!`$G`$20:`$G`$22)"
$ws.Range("C3").Formula = "=STDEV('[2]20231016_Zeta-PotentialData'!`$G`$20:`$G`$22)"

# --- Row 4: DMA S1 (G2:G4) ---------------------------------------------------
$ws.Range("A4").Value = "DMA S1"
$ws.Range("B4").Formula = "=AVERAGE('[1]20231016_Zeta-PotentialData'!`$G`$2:`$G`$4)"
$ws.Range("C4").Formula = "=STDEV('[2]20231016_Zeta-PotentialData'!`$G`$2:`$G`$4)"

# --- Row 5: DIP B1 (G14:G16) -------------------------------------------------
$ws.Range("A5").Value = "DIP B1"
$ws.Range("B5").Formula = "=AVERAGE('[1]20231016_Zeta-PotentialData'!`$G`$14:`$G`$16)"
$ws.Range("C5").Formula = "=STDEV('[2]20231016_Zeta-PotentialData'!`$G`$14:`$G`$16)"

# --- Row 6: DMA B1 (G17:G19) -------------------------------------------------
$ws.Range("A6").Value = "DMA B1"
$ws.Range("B6").Formula = "=AVERAGE('[1]20231016_Zeta-PotentialData'!`$G`$17:`$G`$19)"
$ws.Range("C6").Formula = "=STDEV('[2]20231016_Zeta-PotentialData'!`$G`$17:`$G`$19)"

# --- Row 7: DIP G2 (G11:G13) -------------------------------------------------
$ws.Range("A7").Value = "DIP G2"
$ws.Range("B7").Formula = "=AVERAGE('[1]20231016_Zeta-PotentialData'!`$G`$11:`$G`$13)"
$ws.Range("C7").Formula = "=STDEV('[2]20231016_Zeta-PotentialData'!`$G`$11:`$G`$13)"

# --- Row 8: DMA G1 (G8:G10) --------------------------------------------------
$ws.Range("A8").Value = "DMA G1"
$ws.Range("B8").Formula = "=AVERAGE('[1]20231016_Zeta-PotentialData'!`$G`$8:`$G`$10)"
$ws.Range("C8").Formula = "=STDEV('[2]20231016_Zeta-PotentialData'!`$G`$8:`$G`$10)"

# --- Row 9: DMA G2 (G5:G7) ---------------------------------------------------
$ws.Range("A9").Value = "DMA G2"
$ws.Range("B9").Formula = "=AVERAGE('[1]20231016_Zeta-PotentialData'!`$G`$5:`$G`$7)"
$ws.Range("C9").Formula = "=STDEV('[2]20231016_Zeta-PotentialData'!`$G`$5:`$G`$7)"

# --- Update the active selection on Sheet1 ---------------------------------
$ws.Range("E8").Select()
